$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 40.32060325799205
$ws.Range("G2").Value = 39.67546037221035
$ws.Range("H2").Value = 40.95434293043651
$ws.Range("I2").Value = 0.0007617811570131555
$ws.Range("J2").Value = 0.0007110563723136485
$ws.Range("K2").Value = 0.0008547685285345647
$ws.Range("L2").Value = 0.05766821399896885
$ws.Range("M2").Value = 0.05723096128650838
$ws.Range("N2").Value = 0.05811550414546574

$ws.Range("F3").Value = 0.000014829001671811689280116101
$ws.Range("G3").Value = 0.000000005023164888960513887317
$ws.Range("H3").Value = 0.00004174173713924094934009032
$ws.Range("I3").Value = 0.000012984768622844429657755931
$ws.Range("J3").Value = 0.000000004693368378058987312605
$ws.Range("K3").Value = 0.000036399347969283488840597229
$ws.Range("L3").Value = 0.000015274929953179420833581093
$ws.Range("M3").Value = 0.000000005189443700012504655984
$ws.Range("N3").Value = 0.000042973148430734496973706915

$ws.Range("F4").Value = 40.32061808699373
$ws.Range("G4").Value = 39.67546037723351
$ws.Range("H4").Value = 40.95438467217365
$ws.Range("I4").Value = 0.0007747659256359999
$ws.Range("J4").Value = 0.0007110610656820265
$ws.Range("K4").Value = 0.0008911678765038482
$ws.Range("L4").Value = 0.05768348892892202
$ws.Range("M4").Value = 0.05723096647595208
$ws.Range("N4").Value = 0.05815847729389648
